# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2-155) from 2023-09-10 (45179) to 2023-09-11 (45180).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C155").Value = 45180
